# "Feat (chp3): Correction exo plus"
#
# Slide 1 ("Chapitre 8 ...") title slide, subtitle shape "Sous-titre 2":
# the single run
#     "Chapitre 8: Les chaînes de caractères"
# is split in two runs, with a space added after "caractères":
#     "Chapitre 8: Les chaînes "   (run 1, unchanged)
#     "de caractères "            (run 2, now ends with a trailing space)

$pres = $ppt.ActivePresentation
$slide = $pres.Slides.Item(1)
$shape = $slide.Shapes.Item("Sous-titre 2")
$textRange = $shape.TextFrame.TextRange

# Full original text: "Chapitre 8: Les chaînes de caractères" (37 chars).
# Characters 1-24  -> "Chapitre 8: Les chaînes "  (kept as-is)
# Characters 25-37 -> "de caractères"             (re-typed with trailing space)
$secondPart = $textRange.Characters(25, 13)
$secondPart.Text = "de caractères "

# Keep the explicit font size consistent with the rest of the line (28pt).
$secondPart = $shape.TextFrame.TextRange.Characters(25, 14)
$secondPart.Font.Size = 28
